# Auto-generated Excel COM-interop script to update Zodiark_Profits market data values.
# Applies per-cell numeric updates (and a few cell additions/removals) across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 2281.9
$ws.Range("I115").Value = 1424.4445
$ws.Range("J115").Value = 9999
$ws.Range("K115").Value = 4273.333500000001
$ws.Range("L115").Value = 29997
$ws.Range("M115").Value = -2706.333500000001
$ws.Range("N115").Value = -33131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1982.2787
$ws.Range("I32").Value = 1219.9272
$ws.Range("K32").Value = 1219.9272
$ws.Range("M32").Value = -932.9272000000001

$ws.Range("H74").Value = 1871.3513
$ws.Range("I74").Value = 1680.5
$ws.Range("K74").Value = 1680.5
$ws.Range("M74").Value = -806.5

$ws.Range("H77").Value = 1871.3513
$ws.Range("I77").Value = 1680.5
$ws.Range("K77").Value = 8402.5
$ws.Range("M77").Value = -4034.5

$ws.Range("H94").Value = 71271.82000000001
$ws.Range("J94").Value = 71271.82000000001
$ws.Range("L94").Value = 71271.82000000001
$ws.Range("N94").Value = -73073.82000000001

$ws.Range("H132").Value = 5335.5654
$ws.Range("I132").Value = 4330.4316
$ws.Range("K132").Value = 12991.2948
$ws.Range("M132").Value = -10461.2948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1349.5
$ws.Range("I134").Value = 1349.5
$ws.Range("K134").Value = 4048.5
$ws.Range("M134").Value = -1513.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5813.963
$ws.Range("I31").Value = 2746.1667
$ws.Range("J31").Value = 8268.200000000001
$ws.Range("K31").Value = 2746.1667
$ws.Range("L31").Value = 8268.200000000001
$ws.Range("M31").Value = -2451.1667
$ws.Range("N31").Value = -8858.200000000001

$ws.Range("H34").Value = 5813.963
$ws.Range("I34").Value = 2746.1667
$ws.Range("J34").Value = 8268.200000000001
$ws.Range("K34").Value = 2746.1667
$ws.Range("L34").Value = 8268.200000000001
$ws.Range("M34").Value = -2544.1667
$ws.Range("N34").Value = -8672.200000000001

$ws.Range("H120").Value = 30047.5
$ws.Range("J120").Value = 25799
$ws.Range("L120").Value = 25799
$ws.Range("N120").Value = -33057

$ws.Range("H132").Value = 2223.3845
$ws.Range("I132").Value = 1918.9048
$ws.Range("K132").Value = 5756.7144
$ws.Range("M132").Value = -3226.7144

$ws.Range("H134").Value = 3548.8462
$ws.Range("I134").Value = 3698.55
$ws.Range("K134").Value = 11095.65
$ws.Range("M134").Value = -8560.650000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 824
$ws.Range("I7").Value = 799.6667
$ws.Range("K7").Value = 2399.0001
$ws.Range("M7").Value = -2287.0001

$ws.Range("H39").Value = 3300
$ws.Range("J39").Value = 3766.6667
$ws.Range("L39").Value = 11300.0001
$ws.Range("N39").Value = -11888.0001

$ws.Range("H41").Value = 1791.3334
$ws.Range("J41").Value = 1187.5
$ws.Range("L41").Value = 3562.5
$ws.Range("N41").Value = -4238.5

$ws.Range("H64").Value = 11231.167
$ws.Range("I64").Value = 10479
$ws.Range("K64").Value = 31437
$ws.Range("M64").Value = -31167

$ws.Range("H67").Value = 11231.167
$ws.Range("I67").Value = 10479
$ws.Range("K67").Value = 31437
$ws.Range("M67").Value = -30501

$ws.Range("H68").Value = 693.8333
$ws.Range("I68").Value = 391.33334
$ws.Range("K68").Value = 1174.00002
$ws.Range("M68").Value = -363.0000199999999

$ws.Range("H71").Value = 693.8333
$ws.Range("I71").Value = 391.33334
$ws.Range("K71").Value = 3522.00006
$ws.Range("M71").Value = 533.9999399999997

$ws.Range("H87").Value = 17000.889
$ws.Range("I87").Value = 13584.667
$ws.Range("K87").Value = 40754.001
$ws.Range("M87").Value = -39506.001

$ws.Range("H88").Value = 15999
$ws.Range("J88").Value = 15999
$ws.Range("L88").Value = 47997
$ws.Range("N88").Value = -48853

$ws.Range("H90").Value = 17000.889
$ws.Range("I90").Value = 13584.667
$ws.Range("K90").Value = 122262.003
$ws.Range("M90").Value = -116022.003

$ws.Range("H91").Value = 15999
$ws.Range("J91").Value = 15999
$ws.Range("L91").Value = 47997
$ws.Range("N91").Value = -50961

$ws.Range("H122").Value = 2078.2083
$ws.Range("I122").Value = 1568.7693
$ws.Range("J122").Value = 2680.2727
$ws.Range("K122").Value = 14118.9237
$ws.Range("L122").Value = 24122.4543
$ws.Range("M122").Value = -11668.9237
$ws.Range("N122").Value = -29022.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 29900
$ws.Range("I33").Value = 29900
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 29900
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -29648
$ws.Range("N33").ClearContents()

$ws.Range("H70").Value = 45806.77
$ws.Range("I70").Value = 92171.336
$ws.Range("K70").Value = 92171.336
$ws.Range("M70").Value = -91901.336

$ws.Range("H73").Value = 45806.77
$ws.Range("I73").Value = 92171.336
$ws.Range("K73").Value = 92171.336
$ws.Range("M73").Value = -91235.336

$ws.Range("H126").Value = 4760162.5
$ws.Range("I126").Value = 3365.6667
$ws.Range("K126").Value = 10097.0001
$ws.Range("M126").Value = -7627.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 35721756
$ws.Range("I7").Value = 7656.4
$ws.Range("K7").Value = 7656.4
$ws.Range("M7").Value = -7544.4

$ws.Range("H22").Value = 45291.6
$ws.Range("I22").Value = 127137
$ws.Range("J22").Value = 6776.1177
$ws.Range("K22").Value = 127137
$ws.Range("L22").Value = 6776.1177
$ws.Range("M22").Value = -126842
$ws.Range("N22").Value = -7366.1177

$ws.Range("H27").Value = 45291.6
$ws.Range("I27").Value = 127137
$ws.Range("J27").Value = 6776.1177
$ws.Range("K27").Value = 127137
$ws.Range("L27").Value = 6776.1177
$ws.Range("M27").Value = -127030
$ws.Range("N27").Value = -6990.1177

$ws.Range("H100").Value = 6995.6665
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 6995.6665
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 6995.6665
$ws.Range("N100").Value = -8077.6665
$ws.Range("M100").ClearContents()

$ws.Range("H126").Value = 35721756
$ws.Range("I126").Value = 7656.4
$ws.Range("K126").Value = 22969.2
$ws.Range("M126").Value = -20499.2

$ws.Range("H132").Value = 3804.6667
$ws.Range("I132").Value = 3729.625
$ws.Range("K132").Value = 11188.875
$ws.Range("M132").Value = -8658.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1417.5385
$ws.Range("I132").Value = 1419.5
$ws.Range("K132").Value = 4258.5
$ws.Range("M132").Value = -1728.5

$ws.Range("H136").Value = 5291.385
$ws.Range("I136").Value = 6369.7
$ws.Range("K136").Value = 19109.1
$ws.Range("M136").Value = -16559.1

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
